$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (A4 = "FilesTab") holds the Neo4j query for the Files tab in
# cell B4. The query's RETURN clause is trimmed: the `f.file_type` and
# `demo.breed` columns are dropped.
$newQuery = "`nMATCH (f:file)-->(parent)`nWITH DISTINCT f, parent`nMATCH (f)-[*]->(c:case)<--(demo:demographic)`n MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)`nWHERE diag.stage_of_disease IN ['Unknown']`nWITH DISTINCT f, parent, c, demo, diag, s`nRETURN coalesce(f.file_name, '') AS ``File Name``, `n        coalesce(labels(parent)[0], '') AS ``Association``,`n        coalesce(f.file_description, '') AS ``Description``,`n        coalesce(f.file_format, '') AS ``Format``,`n        coalesce(f.file_size, '') AS ``Size``,`n        coalesce(c.case_id, '') AS ``Case ID``, `n        coalesce(diag.disease_term,'') AS Diagnosis , `n        coalesce(s.clinical_study_designation,'') AS ``Study Code``"

$ws.Range("B4").Value = $newQuery

# Mirror the author's final selection/scroll state: the sheet view ends
# up scrolled to row 4 with B4 (the edited cell) selected.
$ws.Range("B4").Select()
$excel.ActiveWindow.ScrollRow = 4
